$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1537
$ws1.Range("F6").Value = 13481
$ws1.Range("F7").Value = 13315
$ws1.Range("F8").Value = 1029
$ws1.Range("F9").Value = 787
$ws1.Range("F11").Value = 575
$ws1.Range("F13").Value = 7
$ws1.Range("F14").Value = 9
$ws1.Range("F15").Value = 703
$ws1.Range("F16").Value = 2111
$ws1.Range("F17").Value = 32
$ws1.Range("F20").Value = 87
$ws1.Range("F22").Value = 411
$ws1.Range("F23").Value = 300
$ws1.Range("F25").Value = 455
$ws1.Range("F26").Value = 779
$ws1.Range("F27").Value = 39

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 140
$ws2.Range("F8").Value = 658

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 204
$ws3.Range("F3").Value = 71

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 204
$ws4.Range("F4").Value = 1537
$ws4.Range("F8").Value = 13481
$ws4.Range("F9").Value = 13315
$ws4.Range("F10").Value = 1029
$ws4.Range("F11").Value = 787
$ws4.Range("F13").Value = 575
$ws4.Range("F15").Value = 7
$ws4.Range("F16").Value = 9
$ws4.Range("F17").Value = 703
$ws4.Range("F20").Value = 2111
$ws4.Range("F21").Value = 32
$ws4.Range("F24").Value = 87
$ws4.Range("F28").Value = 71
$ws4.Range("F29").Value = 411
$ws4.Range("F30").Value = 300
$ws4.Range("F32").Value = 455
$ws4.Range("F33").Value = 779
$ws4.Range("F34").Value = 140
$ws4.Range("F35").Value = 658
$ws4.Range("F38").Value = 39
